# Refresh the quarterly "Value" series in column B with newly revised source
# data (rows 2-14). Column C holds a YoY % change formula that recalculates
# automatically from these inputs.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$updates = @{
    2  = 3.3439999999999999
    3  = 3.274
    4  = 3.234
    5  = 3.149
    6  = 3.069
    7  = 3.03
    8  = 2.9830000000000001
    9  = 2.9009999999999998
    10 = 2.835
    11 = 2.7730000000000001
    12 = 2.754
    13 = 2.7229999999999999
    14 = 2.6779999999999999
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 2).Value = $updates[$row]
}
